# Figures 2 and 3, report work
#
# Slide 2 ("Figure 2...") caption textbox: resize to span the full slide
# width and bump the body text from 12pt to 14pt.
# Slide 3 ("Figure 3...") caption textbox: resize to span the full slide
# width (font there is already 14pt).

$p = $ppt.ActivePresentation

# --- Slide 2 : "Figure 2." caption textbox (TextBox 4, shape index 1) ---
$s2 = $p.Slides.Item(2)
$tb2 = $s2.Shapes.Item(1)

$tb2.Left = 0
$tb2.Top = 302.1685039370079
$tb2.Width = 960
$tb2.Height = 126.0187421574803

$tb2.TextFrame.TextRange.Font.Size = 14

# --- Slide 3 : "Figure 3." caption textbox (TextBox 4, shape index 1) ---
$s3 = $p.Slides.Item(3)
$tb3 = $s3.Shapes.Item(1)

$tb3.Left = -0.00015748031496062991
$tb3.Top = 350.4717322834646
$tb3.Width = 960.0000917401575
$tb3.Height = 75.12653743307087
